# This script reproduces the edit described by the diff:
#  - Row 4 (old "A=3" entry) is removed entirely, shifting every row below it
#    up by one (so the old row 5 becomes the new row 4, etc.).
#  - The final (now-duplicate) last row is also removed, shrinking the
#    table from A1:D107 down to A1:D105.
#  - The rank numbers in column A for the first two data rows (rows 2 and 3)
#    are cleared out (blanked), matching the diff where A2/A3 become empty
#    inlineStr cells instead of numeric 1/2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire 4th row - this shifts rows 5..107 up to become rows 4..106
$ws.Rows("4").Delete()

# Delete the entire 107th row (now the used range's last row, 106, holds what
# used to be old row 107's data after the shift above) - remove it so the
# table ends at row 105 just like the target dimension A1:D105
$ws.Rows("106").Delete()

# Clear out the rank values in A2 and A3 (they become blank cells)
$ws.Range("A2").Value = ""
$ws.Range("A3").Value = ""
